$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '44.205.99'
$ws.Range('E2').Value = '  +1.35%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.262.76'
$ws.Range('E3').Value = '  +2.96%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '99.33'
$ws.Range('E5').Value = '  +17.62%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '275.70'
$ws.Range('E6').Value = '  +7.18%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.622'
$ws.Range('E7').Value = '  +0.35%  '
$ws.Range('E8').Value = '  -0.10%  '
$ws.Range('E9').Value = '  +7.07%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '48.55'
$ws.Range('E10').Value = '  +8.39%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0948'
$ws.Range('E11').Value = '  +3.41%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '8.25'
$ws.Range('E12').Value = '  +14.03%  '
$ws.Range('E13').Value = '  +0.73%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '15.60'
$ws.Range('E14').Value = '  +8.65%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '2.599.31'
$ws.Range('E15').Value = '  +2.76%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.839'
$ws.Range('E16').Value = '  +7.15%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.249.49'
$ws.Range('E17').Value = '  +0.78%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '44.212.49'
$ws.Range('E18').Value = '  +1.47%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.0000106'
$ws.Range('E19').Value = '  +3.24%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.25'
$ws.Range('E20').Value = '  +5.85%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '71.10'
$ws.Range('E21').Value = '  +1.93%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '10.88'
$ws.Range('E22').Value = '  +21.94%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '2.36'
$ws.Range('E23').Value = '  -0.34%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '235.41'
$ws.Range('E24').Value = '  +1.86%  '
$ws.Range('E25').Value = '  -0.01%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '11.56'
$ws.Range('E26').Value = '  +8.83%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.54'
$ws.Range('E27').Value = '  +14.33%  '
$ws.Range('B28').Value = 'InjectiveProtocol'
$ws.Range('C28').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '40.13'
$ws.Range('E28').Value = '  +3.12%  '
$ws.Range('B29').Value = 'WEMIXToken'
$ws.Range('C29').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '3.38'
$ws.Range('E29').Value = '  -2.99%  '
$ws.Range('E30').Value = '  +0.41%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '173.61'
$ws.Range('E31').Value = '  +0.17%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.0923'
$ws.Range('E32').Value = '  +7.15%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '21.31'
$ws.Range('E33').Value = '  +4.72%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.78'
$ws.Range('E34').Value = '  +9.05%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.125'
$ws.Range('E35').Value = '  +2.23%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.114'
$ws.Range('E36').Value = '  +3.92%  '
$ws.Range('B37').Value = 'RenderToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '4.45'
$ws.Range('E37').Value = '  -0.54%  '
$ws.Range('B38').Value = 'VeChain'
$ws.Range('C38').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0355'
$ws.Range('E38').Value = '  -0.84%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.59'
$ws.Range('E39').Value = '  +25.54%  '
$ws.Range('E40').Value = '  +28.67%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '12.78'
$ws.Range('E41').Value = '  +2.74%  '
$ws.Range('E42').Value = '  +5.55%  '
$ws.Range('B43').Value = 'MultiversX'
$ws.Range('C43').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '63.06'
$ws.Range('E43').Value = '  -0.08%  '
$ws.Range('B44').Value = 'THORChain'
$ws.Range('C44').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '5.50'
$ws.Range('E44').Value = '  +0.74%  '
$ws.Range('E45').Value = '  +5.55%  '
$ws.Range('E46').Value = '  +3.19%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '100.65'
$ws.Range('E47').Value = '  +0.58%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.17'
$ws.Range('E48').Value = '  +5.36%  '
$ws.Range('E49').Value = '  +1.31%  '
$ws.Range('E50').Value = '  +1.02%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.477.48'
$ws.Range('E51').Value = '  +2.53%  '
